$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Índice"
$ws.Range("B1").Value = "Distancia"
$ws.Range("C1").Value = "max"
$ws.Range("D1").Value = "min"
$ws.Range("E1").Value = "Tempo"

# Data rows
$data = @(
    @(0, 1578, 1578, 1578, 0.01085558732350667),
    @(1, 911,  911,  911,  0.01102190812428792),
    @(2, 1233, 1233, 1233, 0.0128873348236084),
    @(3, 1162, 1162, 1162, 0.01283257007598877),
    @(4, 1356, 1356, 1356, 0.0134286642074585),
    @(5, 1484, 1484, 1484, 0.01296807924906413),
    @(6, 1921, 1921, 1921, 0.01290125052134196),
    @(7, 1748, 1748, 1748, 0.01101494630177816),
    @(8, 1645, 1645, 1645, 0.01294527053833008),
    @(9, 1741, 1741, 1741, 0.01291903654734294)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r++
}
